# Apply cryptos list update (prices and volumes) per commit on Mon Jul 15 19:35:27 UTC 2024
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.508.38'
$ws.Range('E2').Value = '  +5.87%  '
$ws.Range('D3').Value = '3.395.51'
$ws.Range('E3').Value = '  +6.46%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Style = 'Normal'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '577.70'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +7.65%  '
$ws.Range('D6').Style = 'Normal'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '155.53'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +7.32%  '
$ws.Range('D7').Style = 'Normal'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').Value = '3.410.98'
$ws.Range('E8').Value = '  +6.74%  '
$ws.Range('D9').Style = 'Normal'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.534'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.22%  '
$ws.Range('D10').Style = 'Normal'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.53'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.92%  '
$ws.Range('D11').Style = 'Normal'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.121'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +7.35%  '
$ws.Range('D12').Style = 'Normal'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.435'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.68%  '
$ws.Range('D13').Value = '3.979.46'
$ws.Range('E13').Value = '  +6.38%  '
$ws.Range('E14').Value = '  +0.27%  '
$ws.Range('D15').Style = 'Normal'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000185'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +7.50%  '
$ws.Range('D16').Style = 'Normal'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '27.14'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +5.51%  '
$ws.Range('D17').Value = '63.535.98'
$ws.Range('E17').Value = '  +5.92%  '
$ws.Range('D18').Value = '3.409.87'
$ws.Range('E18').Value = '  +7.50%  '
$ws.Range('D19').Style = 'Normal'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.39'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.57%  '
$ws.Range('D20').Style = 'Normal'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.07'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +6.45%  '
$ws.Range('D21').Style = 'Normal'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '8.41'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.78%  '
$ws.Range('D22').Style = 'Normal'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '387.68'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +5.13%  '
$ws.Range('D23').Style = 'Normal'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.997'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.32%  '
$ws.Range('D24').Style = 'Normal'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.536'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.76%  '
$ws.Range('D25').Style = 'Normal'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '70.95'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.07%  '
$ws.Range('D26').Style = 'Normal'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.55'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +11.19%  '
$ws.Range('D27').Style = 'Normal'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0000105'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +20.20%  '
$ws.Range('E28').Value = '  +6.78%  '
$ws.Range('E29').Value = '  +1.41%  '
$ws.Range('D30').Style = 'Normal'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.04'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +8.16%  '
$ws.Range('D31').Style = 'Normal'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.55'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +7.67%  '
$ws.Range('B32').Value = 'NEARProtocol'
$ws.Range('C32').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D32').Style = 'Normal'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.70'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +8.18%  '
$ws.Range('B33').Value = 'Fetch.AI'
$ws.Range('C33').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D33').Style = 'Normal'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.34'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +13.08%  '
$ws.Range('D34').Style = 'Normal'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '23.16'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.14%  '
$ws.Range('E35').Value = '  -0.11%  '
$ws.Range('D36').Style = 'Normal'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.72'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.24%  '
$ws.Range('D37').Style = 'Normal'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.49'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +9.57%  '
$ws.Range('D38').Style = 'Normal'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '157.80'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.06%  '
$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D39').Style = 'Normal'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.88'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +11.15%  '
$ws.Range('B40').Value = 'EnergySwap'
$ws.Range('C40').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '27.51'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +4.53%  '
$ws.Range('E41').Value = '  +8.45%  '
$ws.Range('D42').Value = '2.938.63'
$ws.Range('E42').Value = '  +5.47%  '
$ws.Range('D43').Style = 'Normal'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0322'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +4.48%  '
$ws.Range('D44').Style = 'Normal'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.766'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +6.54%  '
$ws.Range('D45').Style = 'Normal'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '41.35'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +3.51%  '
$ws.Range('D46').Style = 'Normal'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.32'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.70%  '
$ws.Range('E47').Value = '  +9.23%  '
$ws.Range('D48').Style = 'Normal'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '22.51'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +9.50%  '
$ws.Range('D49').Value = '3.440.21'
$ws.Range('E49').Value = '  +6.46%  '
$ws.Range('D50').Style = 'Normal'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.33'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +3.06%  '
$ws.Range('E51').Value = '  -2.20%  '
